$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 259, pushing the existing rows 259-273 down to 260-274.
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new price record.
$ws.Cells.Item(259, 1).Value = 5
$ws.Cells.Item(259, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(259, 3).Value = "Maule"
$ws.Cells.Item(259, 4).Value = 44746
$ws.Cells.Item(259, 5).Value = 7
$ws.Cells.Item(259, 6).Value = "Fruta"
$ws.Cells.Item(259, 7).Value = 100101
$ws.Cells.Item(259, 8).Value = "Berries"
$ws.Cells.Item(259, 9).Value = 100101007
$ws.Cells.Item(259, 10).Value = "Kiwi"
$ws.Cells.Item(259, 11).Value = "Hayward"
$ws.Cells.Item(259, 12).Value = "Primera"
$ws.Cells.Item(259, 13).Value = 200
$ws.Cells.Item(259, 14).Value = 5000
$ws.Cells.Item(259, 15).Value = 5000
$ws.Cells.Item(259, 16).Value = 5000
$ws.Cells.Item(259, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(259, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(259, 19).Value = 278
$ws.Cells.Item(259, 20).Value = 18
